# NIT-9012125601.xlsx - "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker "LAURA ELENA SANMARTIN GONZALES" previously had two rows in the
# account-statement table (one for period 2507, one for period 2503). This
# update drops the old 2507 row for her (duplicated with the 2507 row that
# stays for DAVID CHAVEZ RAMOS) and keeps/updates the remaining row, whose
# period moves from 2503 to 2508 - i.e. row 17 (2507/LAURA) is removed and
# the table shifts up, so what was row 18 (2503/LAURA) becomes the new last
# row 17, now showing period 2508. DAVID CHAVEZ RAMOS's own period also
# rolls forward from 2507 to 2508. The overdue amount and period count at
# the top of the statement are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete duplicate row for LAURA ELENA SANMARTIN GONZALES /
# period 2507 - the sheet shifts up and row 18 becomes the new row 17,
# carrying along its own (bottom-bordered "last row") formatting.
$ws.Rows.Item(17).Delete()

# Roll both workers forward to the new period.
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"

# Refresh the summary figures at the top of the statement.
$ws.Range("E11").Value = 120000
$ws.Range("F13").Value = 1
